# Applies the crypto price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.503.58"
$ws.Range("E2").Value = "  +2.05%  "

$ws.Range("D3").Value = "1.576.83"
$ws.Range("E3").Value = "  +0.28%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.84%  "

$ws.Range("D5").Value = "'212.59"
$ws.Range("E5").Value = "  +0.45%  "

$ws.Range("E6").Value = "  -0.18%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.82%  "

$ws.Range("D8").Value = "'46.25"
$ws.Range("E8").Value = "  +6.69%  "

$ws.Range("D9").Value = "'24.05"
$ws.Range("E9").Value = "  +3.69%  "

$ws.Range("E10").Value = "  -0.69%  "

$ws.Range("D12").Value = "'0.0881"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").Value = "1.801.71"
$ws.Range("E13").Value = "  +0.28%  "

$ws.Range("D14").Value = "1.585.45"
$ws.Range("E14").Value = "  +0.97%  "

$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("E16").Value = "  -1.09%  "

$ws.Range("D17").Value = "28.542.81"
$ws.Range("E17").Value = "  +2.35%  "

$ws.Range("D18").Value = "'62.44"
$ws.Range("E18").Value = "  -1.40%  "

$ws.Range("D19").Value = "'229.72"
$ws.Range("E19").Value = "  +0.45%  "

$ws.Range("E20").Value = "  -0.40%  "

$ws.Range("E21").Value = "  -1.29%  "

$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("D23").Value = "'3.94"
$ws.Range("E23").Value = "  -4.08%  "

$ws.Range("E24").Value = "  -1.30%  "

$ws.Range("E25").Value = "  +4.79%  "

$ws.Range("D26").Value = "'151.43"
$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("E28").Value = "  -1.53%  "

$ws.Range("E29").Value = "  -1.85%  "

$ws.Range("E30").Value = "  +0.74%  "

$ws.Range("E31").Value = "  -1.41%  "

$ws.Range("E32").Value = "  -1.72%  "

$ws.Range("E33").Value = "  -0.40%  "

$ws.Range("E34").Value = "  -0.20%  "

$ws.Range("D35").Value = "1.394.45"
$ws.Range("E35").Value = "  -1.24%  "

$ws.Range("E36").Value = "  -1.67%  "

$ws.Range("D37").Value = "'1.02"
$ws.Range("E37").Value = "  -2.73%  "

$ws.Range("E38").Value = "  +2.63%  "

$ws.Range("E39").Value = "  +6.42%  "

$ws.Range("E40").Value = "  -0.56%  "

$ws.Range("D41").Value = "'0.533"
$ws.Range("E41").Value = "  -1.23%  "

$ws.Range("D43").Value = "'0.795"
$ws.Range("E43").Value = "  -1.18%  "

$ws.Range("E44").Value = "  +0.59%  "

$ws.Range("D45").Value = "'1.85"
$ws.Range("E45").Value = "  +2.47%  "

$ws.Range("D46").Value = "'0.980"
$ws.Range("E46").Value = "  +0.89%  "

$ws.Range("D47").Value = "'62.63"
$ws.Range("E47").Value = "  -1.69%  "

$ws.Range("D48").Value = "1.713.67"
$ws.Range("E48").Value = "  +0.43%  "

$ws.Range("D49").Value = "'85.68"
$ws.Range("E49").Value = "  -1.09%  "

$ws.Range("D51").Value = "0.0₆0103"
$ws.Range("E51").Value = "  +5.62%  "
